$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value as plain text, preserving the "no explicit
# cell style" state the original workbook cells were in (cells like
# D2/D6/etc. hold prices as text, e.g. "95.75" or "42.645.87", and must
# stay text rather than being auto-coerced into numbers by Excel).
function Set-TextCell {
    param(
        $Sheet,
        [string]$Address,
        [string]$Text
    )
    $range = $Sheet.Range($Address)
    # Force text interpretation so strings that look numeric (e.g. "9.01")
    # are not silently converted to a Double.
    $range.NumberFormat = "@"
    $range.Value = $Text
    # Restore the default/general style so the cell ends up exactly like
    # the rest of the untouched text cells (no lingering "Text" style).
    $range.Style = "Normal"
}

Set-TextCell $ws "D2" "42.579.47"
Set-TextCell $ws "E2" "  -0.72%  "
Set-TextCell $ws "D3" "2.294.40"
Set-TextCell $ws "E3" "  -0.22%  "
Set-TextCell $ws "E4" "  +0.01%  "
Set-TextCell $ws "D5" "300.67"
Set-TextCell $ws "E5" "  -1.79%  "
Set-TextCell $ws "D6" "95.75"
Set-TextCell $ws "E6" "  -1.46%  "
Set-TextCell $ws "D7" "0.504"
Set-TextCell $ws "E7" "  -1.21%  "
Set-TextCell $ws "E8" "  +0.09%  "
Set-TextCell $ws "D9" "0.494"
Set-TextCell $ws "E9" "  -1.73%  "
Set-TextCell $ws "D10" "34.58"
Set-TextCell $ws "E10" "  -3.18%  "
Set-TextCell $ws "D11" "19.23"
Set-TextCell $ws "E11" "  +5.52%  "
Set-TextCell $ws "D12" "0.0786"
Set-TextCell $ws "E12" "  -0.67%  "
Set-TextCell $ws "E13" "  +0.24%  "
Set-TextCell $ws "D14" "6.77"
Set-TextCell $ws "E14" "  +0.03%  "
Set-TextCell $ws "D15" "2.654.19"
Set-TextCell $ws "E15" "  -0.15%  "
Set-TextCell $ws "D16" "2.298.93"
Set-TextCell $ws "E16" "  -0.20%  "
Set-TextCell $ws "D17" "0.780"
Set-TextCell $ws "E17" "  -0.40%  "
Set-TextCell $ws "D18" "42.569.83"
Set-TextCell $ws "E18" "  -0.60%  "
Set-TextCell $ws "D19" "12.24"
Set-TextCell $ws "E19" "  -4.85%  "
Set-TextCell $ws "E20" "  -1.82%  "
Set-TextCell $ws "E21" "  -0.66%  "
Set-TextCell $ws "D22" "67.50"
Set-TextCell $ws "E22" "  -0.26%  "
Set-TextCell $ws "D23" "235.16"
Set-TextCell $ws "E23" "  -0.63%  "
Set-TextCell $ws "E24" "  +3.15%  "
Set-TextCell $ws "E25" "  +0.12%  "
Set-TextCell $ws "E26" "  -2.32%  "
Set-TextCell $ws "D27" "24.53"
Set-TextCell $ws "E27" "  -3.46%  "
Set-TextCell $ws "B28" "Monero"
Set-TextCell $ws "C28" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell $ws "D28" "164.59"
Set-TextCell $ws "E28" "  -1.07%  "
Set-TextCell $ws "B29" "Toncoin"
Set-TextCell $ws "C29" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextCell $ws "D29" "2.05"
Set-TextCell $ws "E29" "  -0.38%  "
Set-TextCell $ws "D30" "9.01"
Set-TextCell $ws "D31" "32.17"
Set-TextCell $ws "E31" "  -2.56%  "
Set-TextCell $ws "E32" "  -0.02%  "
Set-TextCell $ws "D33" "4.95"
Set-TextCell $ws "E33" "  -0.95%  "
Set-TextCell $ws "D34" "17.46"
Set-TextCell $ws "E34" "  +0.60%  "
Set-TextCell $ws "B35" "WEMIXToken"
Set-TextCell $ws "C35" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextCell $ws "D35" "2.34"
Set-TextCell $ws "E35" "  -2.35%  "
Set-TextCell $ws "B36" "Hedera"
Set-TextCell $ws "C36" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell $ws "D36" "0.0695"
Set-TextCell $ws "E36" "  +0.29%  "
Set-TextCell $ws "B37" "RenderToken"
Set-TextCell $ws "C37" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell $ws "D37" "4.36"
Set-TextCell $ws "E37" "  -10.00%  "
Set-TextCell $ws "D38" "0.0998"
Set-TextCell $ws "E39" "  -0.36%  "
Set-TextCell $ws "E40" "  +0.13%  "
Set-TextCell $ws "E41" "  -1.35%  "
Set-TextCell $ws "E42" "  +7.61%  "
Set-TextCell $ws "D43" "1.958.26"
Set-TextCell $ws "E43" "  -2.57%  "
Set-TextCell $ws "D44" "10.48"
Set-TextCell $ws "E44" "  +4.69%  "
Set-TextCell $ws "E45" "  -1.08%  "
Set-TextCell $ws "E46" "  -2.98%  "
Set-TextCell $ws "E47" "  -1.38%  "
Set-TextCell $ws "E48" "  -0.98%  "
Set-TextCell $ws "D49" "2.523.30"
Set-TextCell $ws "E49" "  -0.09%  "
Set-TextCell $ws "D50" "52.91"
Set-TextCell $ws "E50" "  -2.09%  "
Set-TextCell $ws "D51" "71.66"
Set-TextCell $ws "E51" "  -0.54%  "
